$p = $ppt.ActivePresentation

# Remove the bad last slide ("Pitanja?" / Questions?) that was appended
# at the end of the deck.
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
